# Update the cryptocurrency price/volume table (columns D and E) on Sheet1
# to reflect the latest scrape, per the automated "Updated cryptos list"
# GitHub Actions commit. Only the D (Price) and E (Volume(1h)) text values
# change; everything else (layout, styles, rows 1/13/46, other columns)
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel's
# automatic type inference turn numeric-looking strings (e.g. "1.00",
# "62.802.09") into actual numbers. We temporarily force the cell to the
# Text number format, assign the value, then restore the cell style back
# to "Normal" so the on-disk style/formatting stays identical to the
# original (the cells never had any custom numeric format to begin with).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "62.802.09"
Set-TextValue "E2" "  -1.74%  "
Set-TextValue "E3" "  -4.89%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "591.86"
Set-TextValue "E5" "  -2.46%  "
Set-TextValue "D6" "134.92"
Set-TextValue "E6" "  -5.89%  "
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "D8" "3.168.69"
Set-TextValue "E8" "  -4.97%  "
Set-TextValue "D9" "0.517"
Set-TextValue "E9" "  -0.84%  "
Set-TextValue "E10" "  -6.75%  "
Set-TextValue "D11" "5.25"
Set-TextValue "E11" "  -5.75%  "
Set-TextValue "D12" "0.453"
Set-TextValue "E12" "  -3.80%  "
Set-TextValue "D14" "34.90"
Set-TextValue "E14" "  -1.10%  "
Set-TextValue "D15" "3.694.83"
Set-TextValue "E15" "  -4.94%  "
Set-TextValue "E16" "  -1.12%  "
Set-TextValue "D17" "3.168.26"
Set-TextValue "E17" "  -5.11%  "
Set-TextValue "D18" "62.802.75"
Set-TextValue "E18" "  -1.92%  "
Set-TextValue "D19" "6.54"
Set-TextValue "E19" "  -5.10%  "
Set-TextValue "D20" "460.97"
Set-TextValue "E20" "  -4.74%  "
Set-TextValue "D21" "13.94"
Set-TextValue "E21" "  -1.60%  "
Set-TextValue "D22" "0.695"
Set-TextValue "E22" "  -6.20%  "
Set-TextValue "E23" "  -4.52%  "
Set-TextValue "D24" "13.42"
Set-TextValue "E24" "  -4.47%  "
Set-TextValue "D25" "83.07"
Set-TextValue "E25" "  -2.43%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "E27" "  -0.05%  "
Set-TextValue "D28" "2.67"
Set-TextValue "E28" "  -4.50%  "
Set-TextValue "D29" "6.78"
Set-TextValue "E29" "  -6.18%  "
Set-TextValue "D30" "7.71"
Set-TextValue "E30" "  -7.53%  "
Set-TextValue "E31" "  -6.28%  "
Set-TextValue "E32" "  -6.35%  "
Set-TextValue "E33" "  -4.47%  "
Set-TextValue "E34" "  -7.28%  "
Set-TextValue "E35" "  -7.12%  "
Set-TextValue "E36" "  -5.01%  "
Set-TextValue "D37" "51.24"
Set-TextValue "E37" "  -2.40%  "
Set-TextValue "E38" "  -6.97%  "
Set-TextValue "D39" "0.0388"
Set-TextValue "E39" "  -3.39%  "
Set-TextValue "D40" "403.11"
Set-TextValue "E40" "  -7.45%  "
Set-TextValue "D41" "8.06"
Set-TextValue "E41" "  -4.02%  "
Set-TextValue "E42" "  -6.02%  "
Set-TextValue "E43" "  -5.79%  "
Set-TextValue "D44" "2.785.69"
Set-TextValue "E44" "  -11.34%  "
Set-TextValue "D45" "0.251"
Set-TextValue "E45" "  -6.75%  "
Set-TextValue "E47" "  -7.38%  "
Set-TextValue "D48" "124.68"
Set-TextValue "E48" "  +0.05%  "
Set-TextValue "D49" "25.21"
Set-TextValue "E49" "  -5.17%  "
Set-TextValue "D50" "34.37"
Set-TextValue "E50" "  -7.79%  "
Set-TextValue "E51" "  -2.33%  "